$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 29   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/19/2022  Through  12/25/2022"

# --- Donor cells for style-preserving paste-special (row 14 is never modified by this diff) ---
# C14 = text "0" with style 14 (General/text)
# E14 = text "***.*" with style 14 (General/text)
# I14 = number with style 15 (#,##0 integer format)
# M14 = number with style 16 (#,##0.0 decimal format)

# --- Cells changing to text "0" (style 14) ---
foreach ($addr in @("C15", "C22", "D23", "C26", "D27")) {
    $ws.Range("C14").Copy()
    $ws.Range($addr).PasteSpecial(-4104)
}

# --- Cells changing to text "***.* " (style 14) ---
foreach ($addr in @("E23", "E27")) {
    $ws.Range("E14").Copy()
    $ws.Range($addr).PasteSpecial(-4104)
}

# --- Cells changing to a number with style 15 (#,##0 integer) ---
$ws.Range("I14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Value = 1

# --- Cells changing to a number with style 16 (#,##0.0 decimal) ---
$ws.Range("M14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("M14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Value = -100

# --- Plain value-only updates (style unchanged) ---
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = 50
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -77.777777777777
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 45
$ws.Range("M16").Value = -26.890756302521
$ws.Range("N16").Value = -86.656441717791
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 84
$ws.Range("J17").Value = 67
$ws.Range("K17").Value = 25.373134328358
$ws.Range("L17").Value = 64.705882352941
$ws.Range("M17").Value = 44.827586206896
$ws.Range("N17").Value = -32.8
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 33.333333333333
$ws.Range("I18").Value = 118
$ws.Range("J18").Value = 92
$ws.Range("K18").Value = 28.260869565217
$ws.Range("L18").Value = -7.8125
$ws.Range("M18").Value = -12.592592592592
$ws.Range("N18").Value = -91.411935953420
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 66.666666666666
$ws.Range("F19").Value = 22
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -54.166666666666
$ws.Range("I19").Value = 482
$ws.Range("J19").Value = 371
$ws.Range("K19").Value = 29.919137466307
$ws.Range("L19").Value = 63.389830508474
$ws.Range("M19").Value = 25.848563968668
$ws.Range("N19").Value = -52.135054617676
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 166.666666666667
$ws.Range("I20").Value = 107
$ws.Range("J20").Value = 69
$ws.Range("K20").Value = 55.072463768115
$ws.Range("L20").Value = 64.615384615384
$ws.Range("M20").Value = 8.080808080808
$ws.Range("N20").Value = -96.786786786786
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 7.142857142857
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 80
$ws.Range("H21").Value = -22.5
$ws.Range("I21").Value = 895
$ws.Range("J21").Value = 669
$ws.Range("K21").Value = 33.781763826606
$ws.Range("L21").Value = 48.91846921797
$ws.Range("M21").Value = 12.015018773466
$ws.Range("N21").Value = -86.237121328617
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = -57.142857142857
$ws.Range("J22").Value = 19
$ws.Range("K22").Value = 89.473684210526
$ws.Range("C24").Value = 35
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 16.666666666666
$ws.Range("F24").Value = 161
$ws.Range("G24").Value = 144
$ws.Range("H24").Value = 11.805555555555
$ws.Range("I24").Value = 1769
$ws.Range("J24").Value = 1384
$ws.Range("K24").Value = 27.817919075144
$ws.Range("L24").Value = 63.796296296296
$ws.Range("M24").Value = 80.510204081632
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 33.333333333333
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = -20
$ws.Range("I25").Value = 188
$ws.Range("J25").Value = 181
$ws.Range("K25").Value = 3.867403314917
$ws.Range("L25").Value = 35.251798561151
$ws.Range("M25").Value = -12.558139534883
$ws.Range("J26").Value = 19
$ws.Range("K26").Value = 10.526315789473
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 0
